$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.996.01"
$ws.Range("E2").Value = "  +3.17%  "

$ws.Range("D3").Value = "3.031.84"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.45"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.60"
$ws.Range("E6").Value = "  +8.28%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.028.63"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.95"
$ws.Range("E10").Value = "  +17.07%  "

$ws.Range("E11").Value = "  +4.19%  "

$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("E13").Value = "  +3.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.65"
$ws.Range("E14").Value = "  +5.02%  "

$ws.Range("D16").Value = "3.534.83"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("D18").Value = "62.977.87"
$ws.Range("E18").Value = "  +2.91%  "

$ws.Range("D19").Value = "3.034.48"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.60"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("E22").Value = "  +3.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +4.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.20"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.47"
$ws.Range("E25").Value = "  +12.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  +9.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.42"
$ws.Range("E27").Value = "  +4.68%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  +12.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.50"
$ws.Range("E30").Value = "  +6.76%  "

$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.64"
$ws.Range("E33").Value = "  +2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("E34").Value = "  +3.04%  "

$ws.Range("E36").Value = "  +2.95%  "

$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("E38").Value = "  +12.56%  "

$ws.Range("E39").Value = "  +9.23%  "

$ws.Range("E40").Value = "  +3.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.42"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.09"
$ws.Range("E42").Value = "  +1.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.310"
$ws.Range("E43").Value = "  +17.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.56"
$ws.Range("E44").Value = "  +15.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.06"
$ws.Range("E45").Value = "  +1.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("E46").Value = "  +4.11%  "

$ws.Range("D47").Value = "2.721.79"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.28"
$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.66"
$ws.Range("E49").Value = "  +11.53%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("E51").Value = "  +8.14%  "
